$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Expand the Excel Table to the new size first so the header row
# --- assigns the right column names as we fill them in below.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G4"))

# --- Header row (row 1) ---
$ws.Range("A1").Value = "description"
$ws.Range("B1").Value = "s.name"
$ws.Range("C1").Value = "f.ids.path"
$ws.Range("D1").Value = "l.name"
$ws.Range("E1").Value = "creator.func"
$ws.Range("F1").Value = "read.driver"
$ws.Range("G1").Value = "write.driver"

# --- Row 2 ---
$ws.Range("A2").Value = "A copy of the user-provided AOI spatial file."
$ws.Range("B2").Value = "aoi.orig"
$ws.Range("C2").Value = "dir_0010_vector"
$ws.Range("D2").Value = "AOI_original.gpkg"
$ws.Range("E2").Value = "none"
$ws.Range("F2").Value = "st_read()"
$ws.Range("G2").Value = "st_write()"

# --- Row 3 ---
$ws.Range("A3").Value = "AOI after being snapped to the nearest 100m extent."
$ws.Range("B3").Value = "aoi.snapped"
$ws.Range("C3").Value = "dir_1010_vector"
$ws.Range("D3").Value = "AOI_snapped.gpkg"
$ws.Range("E3").Value = "aoi_snap2()"
$ws.Range("F3").Value = "st_read()"
$ws.Range("G3").Value = "st_write()"

# --- Row 4 (new row, A4 left blank) ---
$ws.Range("B4").Value = "aoi.r.template"
$ws.Range("C4").Value = "dir_1010_vector"
$ws.Range("D4").Value = "AOI_template_*m.tif"
$ws.Range("E4").Value = "create_template()"
$ws.Range("F4").Value = "rast()"
$ws.Range("G4").Value = "writeRaster()"

# --- Column widths (character units, closest achievable to the stored
# --- OOXML widths of 49.85546875, 13.7109375, 17.85546875, 20.42578125,
# --- 17.7109375, 17.7109375, 13.140625) ---
$ws.Columns.Item(1).ColumnWidth = 49.02083333333333
$ws.Columns.Item(2).ColumnWidth = 12.833333333333334
$ws.Columns.Item(3).ColumnWidth = 17.0
$ws.Columns.Item(4).ColumnWidth = 19.666666666666668
$ws.Columns.Item(5).ColumnWidth = 16.833333333333332
$ws.Columns.Item(6).ColumnWidth = 16.833333333333332
$ws.Columns.Item(7).ColumnWidth = 12.333333333333334

# --- Zoom level ---
$excel.ActiveWindow.Zoom = 140

# --- Selection / active cell ---
$ws.Range("E3").Select()
